$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.780611209335849
$ws.Range("D2").Value = 6.554842872365872
$ws.Range("E2").Value = 24.33813070574213
$ws.Range("F2").Value = 41.37386048936481
$ws.Range("G2").Value = 3.592891420298954
$ws.Range("M2").Value = 42.10400818103771
$ws.Range("B3").Value = 7.708837187108705
$ws.Range("D3").Value = 6.525881002224815
$ws.Range("E3").Value = 22.73657062656616
$ws.Range("F3").Value = 40.18859129926754
$ws.Range("G3").Value = 3.606541551222986
$ws.Range("M3").Value = 39.76429407088607
$ws.Range("B4").Value = 7.666334202181588
$ws.Range("D4").Value = 6.521170801178379
$ws.Range("E4").Value = 21.74130068265563
$ws.Range("F4").Value = 39.50198044236849
$ws.Range("G4").Value = 3.615240326429777
$ws.Range("M4").Value = 38.26281918847521
$ws.Range("B5").Value = 7.649425672114651
$ws.Range("D5").Value = 6.522410507095225
$ws.Range("E5").Value = 21.32247199164257
$ws.Range("F5").Value = 39.23289263602974
$ws.Range("G5").Value = 3.618866339924081
$ws.Range("M5").Value = 37.63504653945539
$ws.Range("B6").Value = 7.646643401729598
$ws.Range("D6").Value = 6.522803160040419
$ws.Range("E6").Value = 21.25212790520169
$ws.Range("F6").Value = 39.18886742791262
$ws.Range("G6").Value = 3.619473378886842
$ws.Range("M6").Value = 37.52985829121396
$ws.Range("B7").Value = 7.666104477476813
$ws.Range("D7").Value = 6.521174908151983
$ws.Range("E7").Value = 21.73570575737628
$ws.Range("F7").Value = 39.49830761372419
$ws.Range("G7").Value = 3.615288897627604
$ws.Range("M7").Value = 38.25441663536925
$ws.Range("B8").Value = 7.755549314939166
$ws.Range("D8").Value = 6.54207114965147
$ws.Range("E8").Value = 23.79609060260449
$ws.Range("F8").Value = 40.9568716531921
$ws.Range("G8").Value = 3.597532905890489
$ws.Range("M8").Value = 41.31095435665283
$ws.Range("B9").Value = 7.942564409447954
$ws.Range("D9").Value = 6.691648760314268
$ws.Range("E9").Value = 27.48083627289002
$ws.Range("F9").Value = 44.12573426105983
$ws.Range("G9").Value = 3.565169392868816
$ws.Range("M9").Value = 46.77685254303503
$ws.Range("B10").Value = 8.085917296137495
$ws.Range("D10").Value = 6.873036857174943
$ws.Range("E10").Value = 29.90692789391282
$ws.Range("F10").Value = 46.61699785827755
$ws.Range("G10").Value = 3.542796165870817
$ws.Range("M10").Value = 50.45863817066585
$ws.Range("B11").Value = 8.15217857099201
$ws.Range("D11").Value = 6.97180091917764
$ws.Range("E11").Value = 30.95077145099389
$ws.Range("F11").Value = 47.78051796666346
$ws.Range("G11").Value = 3.532901305956774
$ws.Range("M11").Value = 52.05976347470166
$ws.Range("B12").Value = 8.177400052170675
$ws.Range("D12").Value = 7.011569132233884
$ws.Range("E12").Value = 31.33754869612849
$ws.Range("F12").Value = 48.22505492451754
$ws.Range("G12").Value = 3.529193176021548
$ws.Range("M12").Value = 52.65542008238913
$ws.Range("B13").Value = 8.171962718463305
$ws.Range("D13").Value = 7.002898684755259
$ws.Range("E13").Value = 31.25462594616826
$ws.Range("F13").Value = 48.12914695513044
$ws.Range("G13").Value = 3.529990091109993
$ws.Range("M13").Value = 52.52760902365626
$ws.Range("B14").Value = 8.154251077208682
$ws.Range("D14").Value = 6.975025089459091
$ws.Range("E14").Value = 30.98276171633875
$ws.Range("F14").Value = 47.81701359812107
$ws.Range("G14").Value = 3.532595468364331
$ws.Range("M14").Value = 52.10898218576286
$ws.Range("B15").Value = 8.143418460614688
$ws.Range("D15").Value = 6.958260778765489
$ws.Range("E15").Value = 30.81513248146335
$ws.Range("F15").Value = 47.62632410365198
$ws.Range("G15").Value = 3.53419633891858
$ws.Range("M15").Value = 51.85117278416921
$ws.Range("B16").Value = 8.081606274565937
$ws.Range("D16").Value = 6.866912505907741
$ws.Range("E16").Value = 29.83751503874745
$ws.Range("F16").Value = 46.5415315372103
$ws.Range("G16").Value = 3.543448355331773
$ws.Range("M16").Value = 50.35251170408679
$ws.Range("B17").Value = 8.04394021893312
$ws.Range("D17").Value = 6.815063642538275
$ws.Range("E17").Value = 29.22252942733883
$ws.Range("F17").Value = 45.8834865839932
$ws.Range("G17").Value = 3.549195337273923
$ws.Range("M17").Value = 49.41419533518621
$ws.Range("B18").Value = 8.022375732100512
$ws.Range("D18").Value = 6.786770369646505
$ws.Range("E18").Value = 28.86317162366551
$ws.Range("F18").Value = 45.50786472812523
$ws.Range("G18").Value = 3.552527586640969
$ws.Range("M18").Value = 48.86756018701534
$ws.Range("B19").Value = 8.015092192306215
$ws.Range("D19").Value = 6.777452063346791
$ws.Range("E19").Value = 28.74052869729338
$ws.Range("F19").Value = 45.38119131796785
$ws.Range("G19").Value = 3.553660475312003
$ws.Range("M19").Value = 48.68128917414234
$ws.Range("B20").Value = 8.047939635095927
$ws.Range("D20").Value = 6.820424527830261
$ws.Range("E20").Value = 29.28857815011054
$ws.Range("F20").Value = 45.95324283305757
$ws.Range("G20").Value = 3.548580806958188
$ws.Range("M20").Value = 49.51479952611811
$ws.Range("B21").Value = 8.159450061412615
$ws.Range("D21").Value = 6.983147798034207
$ws.Range("E21").Value = 31.06284484195912
$ws.Range("F21").Value = 47.90859095925764
$ws.Range("G21").Value = 3.531829166970986
$ws.Range("M21").Value = 52.23223240564722
$ws.Range("B22").Value = 8.233075288506656
$ws.Range("D22").Value = 7.103304273995603
$ws.Range("E22").Value = 32.17291149430562
$ws.Range("F22").Value = 49.20931115655926
$ws.Range("G22").Value = 3.521106455150966
$ws.Range("M22").Value = 53.94613227580376
$ws.Range("B23").Value = 8.193718742755095
$ws.Range("D23").Value = 7.037905162548183
$ws.Range("E23").Value = 31.58494670165238
$ws.Range("F23").Value = 48.51313256704435
$ws.Range("G23").Value = 3.526809372994818
$ws.Range("M23").Value = 53.03708015228532
$ws.Range("B24").Value = 8.046131216399223
$ws.Range("D24").Value = 6.817996154284163
$ws.Range("E24").Value = 29.25873557792091
$ws.Range("F24").Value = 45.9216976008098
$ws.Range("G24").Value = 3.5488585480375
$ws.Range("M24").Value = 49.46933874536128
$ws.Range("B25").Value = 7.890843387508639
$ws.Range("D25").Value = 6.638938210829992
$ws.Range("E25").Value = 26.53376720229295
$ws.Range("F25").Value = 43.2379979852612
$ws.Range("G25").Value = 3.573670980035683
$ws.Range("M25").Value = 45.35621729872385
